$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A31").Value = "test"
